# Locate the target paragraph: "Como o artigo de resiliência serve?"
$d = $word.ActiveDocument

$target = $null
foreach ($p in $d.Paragraphs) {
    if ($p.Range.Text -like "Como o artigo de resiliência serve?*") {
        $target = $p
        break
    }
}

if ($target -eq $null) {
    throw "Target paragraph not found"
}

# Replace the paragraph's content (text + trailing paragraph mark) with the four
# paragraphs described in the diff: the original sentence (now on its own,
# unchanged), followed by three new strikethrough task paragraphs. The
# "_GoBack" bookmark moves into the second (new) paragraph.
$r1 = $target.Range.Duplicate
$xml1 = '<pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage"><pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml"><pkg:xmlData><w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:body><w:p><w:pPr><w:rPr><w:strike/></w:rPr></w:pPr><w:r><w:rPr><w:strike/></w:rPr><w:t>Como o artigo de resiliência serve?</w:t></w:r></w:p><w:p><w:pPr><w:rPr><w:strike/></w:rPr></w:pPr><w:bookmarkStart w:id="0" w:name="_GoBack"/><w:bookmarkEnd w:id="0"/><w:proofErr w:type="gramStart"/><w:r><w:rPr><w:strike/></w:rPr><w:t>1</w:t></w:r><w:proofErr w:type="gramEnd"/><w:r><w:rPr><w:strike/></w:rPr><w:t xml:space="preserve"> Ver quais combinações de parâmetro serve para plasticidade média e alta para baixo e alto custo. </w:t></w:r></w:p><w:p><w:pPr><w:rPr><w:strike/></w:rPr></w:pPr><w:proofErr w:type="gramStart"/><w:r><w:rPr><w:strike/></w:rPr><w:t>3</w:t></w:r><w:proofErr w:type="gramEnd"/><w:r><w:rPr><w:strike/></w:rPr><w:t xml:space="preserve"> Salvar output de valores </w:t></w:r><w:proofErr w:type="spellStart"/><w:r><w:rPr><w:strike/></w:rPr><w:t>pré</w:t></w:r><w:proofErr w:type="spellEnd"/><w:r><w:rPr><w:strike/></w:rPr><w:t xml:space="preserve">-perturbação e pós-perturbação, ajustando no </w:t></w:r><w:proofErr w:type="spellStart"/><w:r><w:rPr><w:strike/></w:rPr><w:t>behavioral</w:t></w:r><w:proofErr w:type="spellEnd"/><w:r><w:rPr><w:strike/></w:rPr><w:t xml:space="preserve"> </w:t></w:r><w:proofErr w:type="spellStart"/><w:r><w:rPr><w:strike/></w:rPr><w:t>space</w:t></w:r><w:proofErr w:type="spellEnd"/><w:r><w:rPr><w:strike/></w:rPr><w:t xml:space="preserve">.   500 </w:t></w:r><w:proofErr w:type="spellStart"/><w:r><w:rPr><w:strike/></w:rPr><w:t>ticks</w:t></w:r><w:proofErr w:type="spellEnd"/><w:r><w:rPr><w:strike/></w:rPr><w:t xml:space="preserve"> para salvar </w:t></w:r><w:proofErr w:type="spellStart"/><w:r><w:rPr><w:strike/></w:rPr><w:t>pré</w:t></w:r><w:proofErr w:type="spellEnd"/><w:r><w:rPr><w:strike/></w:rPr><w:t xml:space="preserve">-perturbação, perturbação </w:t></w:r><w:proofErr w:type="spellStart"/><w:r><w:rPr><w:strike/></w:rPr><w:t>tick</w:t></w:r><w:proofErr w:type="spellEnd"/><w:r><w:rPr><w:strike/></w:rPr><w:t xml:space="preserve"> 501 e pós-perturbação 1000 </w:t></w:r><w:proofErr w:type="spellStart"/><w:r><w:rPr><w:strike/></w:rPr><w:t>ticks</w:t></w:r><w:proofErr w:type="spellEnd"/><w:r><w:rPr><w:strike/></w:rPr><w:t xml:space="preserve">. </w:t></w:r></w:p><w:p><w:pPr><w:rPr><w:strike/></w:rPr></w:pPr><w:proofErr w:type="gramStart"/><w:r><w:rPr><w:strike/></w:rPr><w:t>5</w:t></w:r><w:proofErr w:type="gramEnd"/><w:r><w:rPr><w:strike/></w:rPr><w:t xml:space="preserve"> Ajustar </w:t></w:r><w:proofErr w:type="spellStart"/><w:r><w:rPr><w:strike/></w:rPr><w:t>behavioral</w:t></w:r><w:proofErr w:type="spellEnd"/><w:r><w:rPr><w:strike/></w:rPr><w:t xml:space="preserve"> </w:t></w:r><w:proofErr w:type="spellStart"/><w:r><w:rPr><w:strike/></w:rPr><w:t>space</w:t></w:r><w:proofErr w:type="spellEnd"/><w:r><w:rPr><w:strike/></w:rPr><w:t xml:space="preserve"> com perturbação, tentando colocar uma forma de colocar 1 contexto ao menos por vez. </w:t></w:r></w:p></w:body></w:document></pkg:xmlData></pkg:part></pkg:package>'
$r1.InsertXML($xml1) | Out-Null

# Re-find the last of the four paragraphs just inserted (ends with "por vez. ")
# and append a new, empty strikethrough paragraph after it.
$last = $null
foreach ($p in $d.Paragraphs) {
    if ($p.Range.Text -like "*por vez. *") {
        $last = $p
    }
}

if ($last -eq $null) {
    throw "Inserted paragraph not found"
}

$r2 = $last.Range.Duplicate
$r2.Collapse(0)
$xml2 = '<pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage"><pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml"><pkg:xmlData><w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:body><w:p><w:pPr><w:rPr><w:strike/></w:rPr></w:pPr></w:p></w:body></w:document></pkg:xmlData></pkg:part></pkg:package>'
$r2.InsertXML($xml2) | Out-Null
